# Recolor every run that currently uses the dark-gray "111111" text color
# (decimal 1118481 == 0x111111) to pure black ("000000" / 0) throughout the
# document. This mirrors the commit's intent of normalizing the table's
# near-black text color to true black.
#
# We walk every paragraph in the document and only touch the ones whose
# font color is exactly the old value, leaving anything else (including the
# lone run after the table that has no explicit color at all) untouched.
#
# The color is applied to a range that stops one character short of the
# paragraph's end (i.e. excludes the paragraph mark) so that we only ever
# rewrite the run-level <w:color> and never introduce paragraph-mark run
# properties (<w:pPr><w:rPr>) that weren't present in the original file.

$d = $word.ActiveDocument

$oldColor = 1118481   # 0x111111
$newColor = 0         # 0x000000

$paras = $d.Paragraphs
$changed = 0

for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $paras.Item($i)
    $rng = $para.Range

    if ($rng.Font.Color -eq $oldColor) {
        $target = $d.Range($rng.Start, $rng.End - 1)
        $target.Font.Color = $newColor
        $changed = $changed + 1
    }
}

Write-Output "Recolored $changed run(s) from 111111 to 000000"
